$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 193, pushing the existing row 193 (and
# everything below it) down by one. Excel copies formatting from the row
# above on insert, which correctly carries the date style into the new D193.
$ws.Rows.Item(193).Insert()

# Populate the newly inserted row 193 with the new data point.
$ws.Cells.Item(193, 1).Value2 = 10
$ws.Cells.Item(193, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(193, 3).Value2 = "La Araucanía"
$ws.Cells.Item(193, 4).Value2 = 44468
$ws.Cells.Item(193, 5).Value2 = 9
$ws.Cells.Item(193, 6).Value2 = 100112023
$ws.Cells.Item(193, 7).Value2 = "Brócoli"
$ws.Cells.Item(193, 8).Value2 = "Sin especificar"
$ws.Cells.Item(193, 9).Value2 = "Primera"
$ws.Cells.Item(193, 10).Value2 = 1500
$ws.Cells.Item(193, 11).Value2 = 900
$ws.Cells.Item(193, 12).Value2 = 900
$ws.Cells.Item(193, 13).Value2 = 900
$ws.Cells.Item(193, 14).Value2 = "$/unidad"
$ws.Cells.Item(193, 15).Value2 = "Región de O'Higgins"
$ws.Cells.Item(193, 16).Value2 = 900
$ws.Cells.Item(193, 17).Value2 = 1
$ws.Cells.Item(193, 18).Value2 = "Hortaliza"
